# Final table fixed with CCs instead of speedup
#
# The "Speed-up comp." table (4th table in the document) originally held
# speed-up ratios (e.g. "1", "1.06", "1.17", "1.27"/"1.28"). The commit
# replaces those with the actual clock-cycle counts (CCs) for each
# configuration.

$d = $word.ActiveDocument
$t = $d.Tables.Item(4)

# Returns a plain Document.Range (not the transient Cell.Range) covering
# the given cell, so that Find/Replace and other mutations reliably persist.
function Get-CellDocRange($table, $row, $col) {
    $cell = $table.Cell($row, $col)
    $cr = $cell.Range
    return $d.Range($cr.Start, $cr.End)
}

# Replace the whole (single-run) text of a cell with new text, scoped
# strictly to that cell (Wrap:=wdFindStop, Replace:=wdReplaceOne so the
# match can't leak into other cells/paragraphs that share the same text).
function Replace-CellText($table, $row, $col, $oldText, $newText) {
    $rng = Get-CellDocRange $table $row $col
    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, `
                       $true, 0, $false, $newText, 1) | Out-Null
}

# Replace a cell's text with two runs of identical character formatting
# (e.g. "1.17" -> "1" + "960", rendering as "1960"). First the old text is
# replaced by the first run's text, then the second run's text is typed
# right before the end-of-cell mark, and finally its formatting is nudged
# and restored so it stays a distinct run instead of re-merging with the
# first one.
function Split-CellText($table, $row, $col, $oldText, $firstPart, $secondPart) {
    $rng = Get-CellDocRange $table $row $col
    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, `
                       $true, 0, $false, $firstPart, 1) | Out-Null

    $cell = $table.Cell($row, $col)
    $cr = $cell.Range
    $insPos = $cr.End - 1
    $insRng = $d.Range($insPos, $insPos)
    $insRng.InsertAfter($secondPart)

    $newLen = $secondPart.Length
    $origSize = $d.Range($insPos, $insPos + $newLen).Font.Size
    $d.Range($insPos, $insPos + $newLen).Font.Size = $origSize + 4
    $d.Range($insPos, $insPos + $newLen).Font.Size = $origSize
}

# Row 2: "By hand"
Replace-CellText $t 2 2 "1" "2240"
Replace-CellText $t 2 3 "1.06" "2120"
Split-CellText   $t 2 4 "1.17" "1" "960"
Split-CellText   $t 2 5 "1.27" "1" "760"

# Row 3: "By simulation"
Replace-CellText $t 3 2 "1" "2210"
Replace-CellText $t 3 3 "1.06" "2090"
Split-CellText   $t 3 4 "1.17" "1" "890"
Split-CellText   $t 3 5 "1.28" "1" "730"
